# Weekly update: insert a new record as the first data row of the
# "Arveja Verde" (Macroferia Regional de Talca) table, pushing the
# existing rows 78-97 down to 79-98.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 78 (the first data row) so every existing
# record shifts down by one row (78->79, ..., 97->98).
$ws.Rows("78:78").Insert()

# Populate the newly inserted row 78 with the new week's record. Columns
# that are constant across every record in this table (A, B, C, E, F, G,
# H, I, N, Q, R) are copied straight from the table's pattern; the ones
# that vary per record (D, J, K, L, M, O, P) carry the new values.
$ws.Range("A78").Value = 5
$ws.Range("B78").Value = "Macroferia Regional de Talca"
$ws.Range("C78").Value = "Maule"
$ws.Range("D78").Value = 44641
$ws.Range("E78").Value = 7
$ws.Range("F78").Value = 100112022
$ws.Range("G78").Value = "Arveja Verde"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 300
$ws.Range("K78").Value = 23000
$ws.Range("L78").Value = 23000
$ws.Range("M78").Value = 23000
$ws.Range("N78").Value = "`$/saco 25 kilos"
$ws.Range("O78").Value = "Carahue"
$ws.Range("P78").Value = 920
$ws.Range("Q78").Value = 25
$ws.Range("R78").Value = "Hortaliza"

# Keep the date-formatted column's number format consistent with the
# rest of column D (the Insert() already carried the style down, but
# make sure explicitly in case the host didn't propagate it).
$ws.Range("D78").NumberFormat = $ws.Range("D79").NumberFormat
